$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure Price column (D) cells are treated as text so values like "1.00" are not
# silently converted to numbers by Excel's automatic type detection.

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '92.472.09'
$ws.Range("E2").Value = '  +6.76%  '

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '3.303.80'
$ws.Range("E3").Value = '  +1.22%  '

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '1.00'
$ws.Range("E4").Value = '  -0.04%  '

$ws.Range("E5").Value = '  +2.21%  '

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '627.81'
$ws.Range("E6").Value = '  +0.35%  '

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.421'
$ws.Range("E7").Value = '  +13.79%  '

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.711'
$ws.Range("E8").Value = '  +3.95%  '

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.998'
$ws.Range("E9").Value = '  -0.13%  '

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '3.292.59'
$ws.Range("E10").Value = '  +1.00%  '

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.594'
$ws.Range("E11").Value = '  +4.17%  '

$ws.Range("E12").Value = '  +5.39%  '

$ws.Range("E13").Value = '  +1.02%  '

$ws.Range("E14").Value = '  +1.98%  '

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '3.905.86'
$ws.Range("E15").Value = '  +0.92%  '

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '91.826.94'
$ws.Range("E16").Value = '  +6.23%  '

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '5.40'
$ws.Range("E17").Value = '  +1.94%  '

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '3.313.68'
$ws.Range("E18").Value = '  +1.34%  '

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '3.35'
$ws.Range("E19").Value = '  +10.44%  '

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '14.19'
$ws.Range("E20").Value = '  +1.63%  '

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '438.66'
$ws.Range("E21").Value = '  +1.97%  '

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '9.02'
$ws.Range("E22").Value = '  +2.09%  '

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '5.35'
$ws.Range("E23").Value = '  +0.92%  '

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '0.0000190'
$ws.Range("E24").Value = '  +48.63%  '

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '5.48'
$ws.Range("E25").Value = '  +8.10%  '

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '12.47'
$ws.Range("E26").Value = '  -0.45%  '

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '3.478.53'
$ws.Range("E27").Value = '  +1.11%  '

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '76.64'
$ws.Range("E28").Value = '  +0.89%  '

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '0.999'
$ws.Range("E29").Value = '  +0.01%  '

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '0.183'
$ws.Range("E30").Value = '  +5.38%  '

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '0.999'
$ws.Range("E31").Value = '  -0.08%  '

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '8.77'
$ws.Range("E32").Value = '  +0.22%  '

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '561.99'
$ws.Range("E33").Value = '  +4.11%  '

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '7.47'
$ws.Range("E34").Value = '  +7.50%  '

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '3.74'
$ws.Range("E35").Value = '  +29.02%  '

$ws.Range("E36").Value = '  -4.41%  '

$ws.Range("E37").Value = '  +0.00%  '

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '22.73'
$ws.Range("E38").Value = '  +1.80%  '

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.133'
$ws.Range("E39").Value = '  -2.28%  '

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '22.41'
$ws.Range("E40").Value = '  +3.90%  '

$ws.Range("E41").Value = '  -0.20%  '

$ws.Range("E42").Value = '  +1.57%  '

$ws.Range("E43").Value = '  +1.08%  '

$ws.Range("E44").Value = '  +0.16%  '

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '182.98'
$ws.Range("E45").Value = '  +2.65%  '

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '148.43'
$ws.Range("E46").Value = '  -4.86%  '

$ws.Range("E47").Value = '  +7.30%  '

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '43.93'
$ws.Range("E48").Value = '  -1.23%  '

$ws.Range("E49").Value = '  +0.61%  '

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '25.62'
$ws.Range("E50").Value = '  +6.93%  '

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '0.633'
$ws.Range("E51").Value = '  +1.43%  '
